$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: status -> "Concluída", comment -> "Aprovação adquirida."
$ws.Range("H22").Value = "Concluída"
$ws.Range("I22").Value = "Aprovação adquirida."

# Row 23: Data Fim -> 23/11/2017 (serial 43062), status -> "Concluída",
# comment -> "Todas as atividades previstas até o momento foram concluídas."
$ws.Range("F23").Value = "11/23/2017"
$ws.Range("H23").Value = "Concluída"
$ws.Range("I23").Value = "Todas as atividades previstas até o momento foram concluídas."

# Row 24: status -> "Concluída", comment -> "Todos confirmaram comprometimento e aprovaram os requisitos."
$ws.Range("H24").Value = "Concluída"
$ws.Range("I24").Value = "Todos confirmaram comprometimento e aprovaram os requisitos."

# Update selection to I24
$ws.Range("I24").Select()
